# The workbook gained one new data row (a new weekly "Ají" price record),
# which was inserted right before the existing row 63. This shifts every
# row from 63..128 down by one (to 64..129) and makes the sheet's
# dimension grow from A1:R128 to A1:R129.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 63, pushing old rows 63-128 to 64-129.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new record's data.
$ws.Range("A63").Value = 11
$ws.Range("B63").Value = "Vega Monumental Concepción"
$ws.Range("C63").Value = "Bíobío"
$ws.Range("D63").Value = 44803
$ws.Range("E63").Value = 8
$ws.Range("F63").Value = 100112021
$ws.Range("G63").Value = "Ají"
$ws.Range("H63").Value = "Inferno"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 18
$ws.Range("K63").Value = 16000
$ws.Range("L63").Value = 17000
$ws.Range("M63").Value = 16556
$ws.Range("N63").Value = "$/caja 12 kilos"
$ws.Range("O63").Value = "Región de Arica y Parinacota"
$ws.Range("P63").Value = 1380
$ws.Range("Q63").Value = 12
$ws.Range("R63").Value = "Hortaliza"
